# edit.ps1 - applies the diff:
#   1. Insert a new "Andre Agassi Tennis" list item before "Barbie: Super Model"
#   2. Move <w:lastRenderedPageBreak/> from the "Carrier Aces" run to the
#      "Captain Novalin" run
#   3. Move <w:lastRenderedPageBreak/> from the "Super Star Wars: Return of
#      the Jedi" run to the "Super " run

$d = $word.ActiveDocument

$pkgOpen = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# ---------------------------------------------------------------------------
# 1. Add the missing "Andre Agassi Tennis" entry right before
#    "Barbie: Super Model" (same list formatting as its neighbours).
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Barbie: Super Model", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.InsertParagraphBefore()

# Re-find "Barbie: Super Model" (its paragraph shifted down by one) and grab
# the newly-created empty paragraph immediately before it.
$rng2 = $d.Content
$null = $rng2.Find.Execute("Barbie: Super Model", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$barbiePara = $rng2.Paragraphs(1)
$newPara = $barbiePara.Previous()
$newPara.Range.Text = "Andre Agassi Tennis"

# ---------------------------------------------------------------------------
# 2. Move the lastRenderedPageBreak marker from "Carrier Aces" to
#    "Captain Novalin".
# ---------------------------------------------------------------------------

# 2a. Add it to the "Captain Novalin" run.
$rng = $d.Content
$null = $rng.Find.Execute("Captain Novalin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$xml = $pkgOpen + "<w:p w14:paraId=`"69843B0A`" w14:textId=`"1C4E318A`" w:rsidR=`"007B4C1B`" w:rsidRDefault=`"007B4C1B`" w:rsidP=`"00434A74`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Captain Novalin</w:t></w:r></w:p>" + $pkgClose
$null = $rng.InsertXML($xml)

# 2b. Remove it from the "Carrier Aces" run.
$rng = $d.Content
$null = $rng.Find.Execute("Carrier Aces", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$xml = $pkgOpen + "<w:p w14:paraId=`"28F0CEAB`" w14:textId=`"49EE2956`" w:rsidR=`"00434A74`" w:rsidRDefault=`"00434A74`" w:rsidP=`"00434A74`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>Carrier Aces</w:t></w:r></w:p>" + $pkgClose
$null = $rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3. Move the lastRenderedPageBreak marker from "Super Star Wars: Return of
#    the Jedi" to the "Super " run (the first run of the "Super Star Wars"
#    paragraph, found positionally as the paragraph right after "Super
#    Hockey" since the plain text "Super " is not unique in this document).
# ---------------------------------------------------------------------------

# 3a. Add it to the "Super " run.
$rng = $d.Content
$null = $rng.Find.Execute("Super Hockey", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hockeyPara = $rng.Paragraphs(1)
$superStarWarsPara = $hockeyPara.Next()
$target = $superStarWarsPara.Range
$target.Collapse(1)
$xml = $pkgOpen + "<w:p w14:paraId=`"1563469D`" w14:textId=`"45A8D2BD`" w:rsidR=`"00A631F0`" w:rsidRDefault=`"00A631F0`" w:rsidP=`"00E23A53`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">Super </w:t></w:r><w:r w:rsidR=`"00B96B5C`"><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>Star Wars</w:t></w:r></w:p>" + $pkgClose
$null = $target.InsertXML($xml)

# 3b. Remove it from the "Super Star Wars: Return of the Jedi" run.
$rng = $d.Content
$null = $rng.Find.Execute("Super Star Wars: Return of the Jedi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$xml = $pkgOpen + "<w:p w14:paraId=`"71687C3F`" w14:textId=`"486F8B20`" w:rsidR=`"00B96B5C`" w:rsidRDefault=`"00B96B5C`" w:rsidP=`"00E23A53`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:t>Super Star Wars: Return of the Jedi</w:t></w:r></w:p>" + $pkgClose
$null = $rng.InsertXML($xml)

Write-Host "Done."
